# The workbook gained one new data row: a new weekly price observation was
# inserted as row 5 (Región del Maule / "Americana (o)"), pushing every
# existing data row (old rows 5..60) down by one (new rows 6..61).
#
# Insert a blank row at row 5 - Excel shifts rows 5..60 down to 6..61
# (carrying the old row 5's style, e.g. the date number format on column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new observation.
$ws.Cells.Item(5, 1).Value  = 11
$ws.Cells.Item(5, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value  = "Bíobío"
$ws.Cells.Item(5, 4).Value  = 44530
$ws.Cells.Item(5, 5).Value  = 8
$ws.Cells.Item(5, 6).Value  = 100112021
$ws.Cells.Item(5, 7).Value  = "Ají"
$ws.Cells.Item(5, 8).Value  = "Americana (o)"
$ws.Cells.Item(5, 9).Value  = "Primera"
$ws.Cells.Item(5, 10).Value = 270
$ws.Cells.Item(5, 11).Value = 17000
$ws.Cells.Item(5, 12).Value = 18000
$ws.Cells.Item(5, 13).Value = 17556
$ws.Cells.Item(5, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(5, 15).Value = "Región del Maule"
$ws.Cells.Item(5, 16).Value = 1170
$ws.Cells.Item(5, 17).Value = 15
$ws.Cells.Item(5, 18).Value = "Hortaliza"
